# Move the K_KREIS row from the bottom (row 12) of the category table to just
# below the header block (row 6), shifting K_LAENDER..K_URBAN down by one row,
# and fix the English translation for K_KREIS from "?" to "County".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contents for rows 6 through 12 (column A = code, B = German, C = English)
$data = @(
    @("K_KREIS",   "Kreis",                  "County"),
    @("K_LAENDER", "Bundesland",             "Federal state"),
    @("K_PM",      "Feinstaub",              "Fine particulate matter"),
    @("K_SEA",     "Meer",                   "Sea"),
    @("K_SERIES",  "Zeitreihe",              "Time series"),
    @("K_SEX",     "Geschlecht",             "Sex"),
    @("K_URBAN",   "Verstädterungsgrad",     "Degree of urbanisation")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
